# Applies: re-sort the "Estado de Cuenta" detail table (rows 16-39) so that
# all records for MARIA EUGENIA GUISADO DE ROMERO (doc 22853823) come first,
# ordered by Periodo Mora descending (2108 -> 2005), followed by all records
# for MILDRETH DEL CARMEN MEZA HENAO (doc 22854431), also ordered by
# Periodo Mora descending (2012 -> 2005). Only columns C (N° Doc Trabajador),
# D (Nombre Trabajador), E (Periodo Mora) and F (Valor Mora) are touched;
# column B (Tipo Doc) and G (Salario Basico) are constant across all rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$firstRow = 16
$lastRow = 39

# Capture the existing (pre-sort) records for the data rows.
$records = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $records += , @{
        B = $ws.Range("B$r").Value2
        C = $ws.Range("C$r").Value2
        D = $ws.Range("D$r").Value2
        E = $ws.Range("E$r").Value2
        F = $ws.Range("F$r").Value2
        G = $ws.Range("G$r").Value2
    }
}

# Sort key: MARIA EUGENIA GUISADO DE ROMERO (22853823) group first, then
# MILDRETH DEL CARMEN MEZA HENAO (22854431) group; within each group,
# Periodo Mora (E) descending. Build a single composite numeric sort key
# per record (since Sort-Object -Property with multiple keys is unreliable
# in this host) and sort ascending on that one key.
$personOrder = @{ "22853823" = 0; "22854431" = 1 }
for ($i = 0; $i -lt $records.Count; $i++) {
    $grp = $personOrder[[string]$records[$i].C]
    $records[$i].SortKey = ($grp * 100000) - [int]$records[$i].E
}

$sorted = $records | Sort-Object -Property SortKey

# Write the re-sorted records back into rows 16-39.
$i = 0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rec = $sorted[$i]
    $ws.Range("B$r").Value2 = $rec.B
    $ws.Range("C$r").Value2 = $rec.C
    $ws.Range("D$r").Value2 = $rec.D
    $ws.Range("E$r").Value2 = $rec.E
    $ws.Range("F$r").Value2 = $rec.F
    $ws.Range("G$r").Value2 = $rec.G
    $i++
}
